# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview zh-cn/de-de status cells, and the per-language Status
#   column on the "zh-cn" / "de-de" detail sheets).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
#   move forward a few seconds to record the new handoff.
# - The Status columns grow wider (content got longer), matching Excel's
#   column auto-fit after the text change.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"
$newWidth  = 16.33   # rounds to the same stored column width Excel uses for the new, longer status text

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-20 13:00:35"

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# ---- zh-cn detail sheet ---------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-20 13:00:31"
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# ---- de-de detail sheet ----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-20 13:00:35"
$dede.Columns.Item(3).ColumnWidth = $newWidth
